$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 260.86667
$ws.Range("I9").Value = 279.14285
$ws.Range("J9").Value = 5
$ws.Range("K9").Value = 279.14285
$ws.Range("L9").Value = 5
$ws.Range("M9").Value = -110.14285
$ws.Range("N9").Value = -343
$ws.Range("H40").Value = 5041.143
$ws.Range("J40").Value = 5464.6665
$ws.Range("L40").Value = 5464.6665
$ws.Range("N40").Value = -5814.6665
$ws.Range("H51").Value = 48446
$ws.Range("I51").Value = 7493.4
$ws.Range("J51").Value = 77697.86
$ws.Range("K51").Value = 7493.4
$ws.Range("L51").Value = 77697.86
$ws.Range("M51").Value = -7009.4
$ws.Range("N51").Value = -78665.86
$ws.Range("H74").Value = 4916.1665
$ws.Range("I74").Value = 4916.1665
$ws.Range("K74").Value = 4916.1665
$ws.Range("M74").Value = -3980.1665
$ws.Range("H76").Value = 4749.3335
$ws.Range("I76").Value = 3749.75
$ws.Range("J76").Value = 5249.125
$ws.Range("K76").Value = 3749.75
$ws.Range("L76").Value = 5249.125
$ws.Range("M76").Value = -3434.75
$ws.Range("N76").Value = -5879.125
$ws.Range("H77").Value = 4916.1665
$ws.Range("I77").Value = 4916.1665
$ws.Range("K77").Value = 24580.8325
$ws.Range("M77").Value = -19900.8325
$ws.Range("H79").Value = 4749.3335
$ws.Range("I79").Value = 3749.75
$ws.Range("J79").Value = 5249.125
$ws.Range("K79").Value = 3749.75
$ws.Range("L79").Value = 5249.125
$ws.Range("M79").Value = -2657.75
$ws.Range("N79").Value = -7433.125
$ws.Range("H80").Value = 67365.77
$ws.Range("J80").Value = 4457.2144
$ws.Range("L80").Value = 13371.6432
$ws.Range("N80").Value = -15367.6432
$ws.Range("H83").Value = 67365.77
$ws.Range("J83").Value = 4457.2144
$ws.Range("L83").Value = 40114.9296
$ws.Range("N83").Value = -50098.9296
$ws.Range("H98").Value = 21078.861
$ws.Range("I98").Value = 24119.88
$ws.Range("J98").Value = 2072.5
$ws.Range("K98").Value = 24119.88
$ws.Range("L98").Value = 2072.5
$ws.Range("M98").Value = -22621.88
$ws.Range("N98").Value = -5068.5
$ws.Range("H106").Value = 6180171.5
$ws.Range("J106").Value = 9372.5
$ws.Range("L106").Value = 9372.5
$ws.Range("N106").Value = -10634.5
$ws.Range("H111").Value = 1259.8182
$ws.Range("I111").Value = 1163
$ws.Range("J111").Value = 1376
$ws.Range("K111").Value = 3489
$ws.Range("L111").Value = 4128
$ws.Range("M111").Value = -422
$ws.Range("N111").Value = -10262
$ws.Range("H122").Value = 21078.861
$ws.Range("I122").Value = 24119.88
$ws.Range("J122").Value = 2072.5
$ws.Range("K122").Value = 72359.64
$ws.Range("L122").Value = 6217.5
$ws.Range("M122").Value = -69909.64
$ws.Range("N122").Value = -11117.5
$ws.Range("H127").Value = 1533.3334
$ws.Range("J127").Value = 1800
$ws.Range("L127").Value = 5400
$ws.Range("N127").Value = -15320
$ws.Range("H137").Value = 7523.3335
$ws.Range("I137").Value = 10174.064
$ws.Range("K137").Value = 30522.192
$ws.Range("M137").Value = -27972.192
$ws.Range("H138").Value = 4521.718
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").Value = $null
$ws.Range("H141").Value = 7962.926
$ws.Range("I141").Value = 8026.087
$ws.Range("J141").Value = 7599.75
$ws.Range("K141").Value = 24078.261
$ws.Range("L141").Value = 22799.25
$ws.Range("M141").Value = -18898.261
$ws.Range("N141").Value = -33159.25

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9217.896000000001
$ws.Range("I32").Value = 8641.537
$ws.Range("J32").Value = 16998.75
$ws.Range("K32").Value = 8641.537
$ws.Range("L32").Value = 16998.75
$ws.Range("M32").Value = -8354.537
$ws.Range("N32").Value = -17572.75
$ws.Range("H61").Value = 11550.131
$ws.Range("I61").Value = 17331
$ws.Range("J61").Value = 6251
$ws.Range("K61").Value = 17331
$ws.Range("L61").Value = 6251
$ws.Range("M61").Value = -17119
$ws.Range("N61").Value = -6675
$ws.Range("H132").Value = 2787.9
$ws.Range("I132").Value = 866.1111
$ws.Range("K132").Value = 2598.3333
$ws.Range("M132").Value = -68.33329999999978
$ws.Range("H136").Value = 11550.131
$ws.Range("I136").Value = 17331
$ws.Range("J136").Value = 6251
$ws.Range("K136").Value = 51993
$ws.Range("L136").Value = 18753
$ws.Range("M136").Value = -49443
$ws.Range("N136").Value = -23853

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2140.65
$ws.Range("I20").Value = 1479.6
$ws.Range("K20").Value = 1479.6
$ws.Range("M20").Value = -1232.6
$ws.Range("H22").Value = 197.5
$ws.Range("I22").Value = 150
$ws.Range("K22").Value = 150
$ws.Range("M22").Value = 23
$ws.Range("H50").Value = 69665.336
$ws.Range("J50").Value = 69998
$ws.Range("L50").Value = 69998
$ws.Range("N50").Value = -71146
$ws.Range("H52").Value = 49966.332
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 49966.332
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 49966.332
$ws.Range("M52").Value = $null
$ws.Range("N52").Value = -50492.332
$ws.Range("H55").Value = 71998
$ws.Range("J55").Value = 73996
$ws.Range("L55").Value = 73996
$ws.Range("N55").Value = -74542
$ws.Range("H64").Value = 7233.0835
$ws.Range("J64").Value = 1679.4
$ws.Range("L64").Value = 1679.4
$ws.Range("N64").Value = -2129.4
$ws.Range("H67").Value = 7233.0835
$ws.Range("J67").Value = 1679.4
$ws.Range("L67").Value = 1679.4
$ws.Range("N67").Value = -3239.4
$ws.Range("H70").Value = 204930
$ws.Range("J70").Value = 204930
$ws.Range("L70").Value = 204930
$ws.Range("N70").Value = -205516
$ws.Range("H73").Value = 204930
$ws.Range("J73").Value = 204930
$ws.Range("L73").Value = 204930
$ws.Range("N73").Value = -206958
$ws.Range("H99").Value = 20979.6
$ws.Range("J99").Value = 7514.1665
$ws.Range("L99").Value = 7514.1665
$ws.Range("N99").Value = -10510.1665
$ws.Range("H121").Value = 49966.332
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 49966.332
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 49966.332
$ws.Range("M121").Value = $null
$ws.Range("N121").Value = -53460.332
$ws.Range("H129").Value = 50709
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").Value = $null
$ws.Range("H134").Value = 6050.1035
$ws.Range("I134").Value = 6588.9
$ws.Range("K134").Value = 19766.7
$ws.Range("M134").Value = -17231.7

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 241.66667
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 241.66667
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 241.66667
$ws.Range("M22").Value = $null
$ws.Range("N22").Value = -941.6666700000001
$ws.Range("H31").Value = 12519.117
$ws.Range("I31").Value = 16297.363
$ws.Range("J31").Value = 5592.3335
$ws.Range("K31").Value = 16297.363
$ws.Range("L31").Value = 5592.3335
$ws.Range("M31").Value = -16002.363
$ws.Range("N31").Value = -6182.3335
$ws.Range("H34").Value = 12519.117
$ws.Range("I34").Value = 16297.363
$ws.Range("J34").Value = 5592.3335
$ws.Range("K34").Value = 16297.363
$ws.Range("L34").Value = 5592.3335
$ws.Range("M34").Value = -16095.363
$ws.Range("N34").Value = -5996.3335
$ws.Range("H58").Value = 2920.1365
$ws.Range("J58").Value = 3679.182
$ws.Range("L58").Value = 3679.182
$ws.Range("N58").Value = -4085.182
$ws.Range("H68").Value = 37167.855
$ws.Range("J68").Value = 37167.855
$ws.Range("L68").Value = 37167.855
$ws.Range("N68").Value = -38665.855
$ws.Range("H71").Value = 37167.855
$ws.Range("J71").Value = 37167.855
$ws.Range("L71").Value = 111503.565
$ws.Range("N71").Value = -118991.565
$ws.Range("H99").Value = 242766.47
$ws.Range("J99").Value = 6212.5
$ws.Range("L99").Value = 6212.5
$ws.Range("N99").Value = -9208.5
$ws.Range("H105").Value = 7608.625
$ws.Range("I105").Value = 10386.818
$ws.Range("K105").Value = 10386.818
$ws.Range("M105").Value = -8639.817999999999
$ws.Range("H126").Value = 242766.47
$ws.Range("J126").Value = 6212.5
$ws.Range("L126").Value = 18637.5
$ws.Range("N126").Value = -23577.5
$ws.Range("H134").Value = 5443
$ws.Range("I134").Value = 3892.889
$ws.Range("K134").Value = 11678.667
$ws.Range("M134").Value = -9143.667000000001
$ws.Range("H136").Value = 2920.1365
$ws.Range("J136").Value = 3679.182
$ws.Range("L136").Value = 11037.546
$ws.Range("N136").Value = -16137.546

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 5500
$ws.Range("J54").Value = 10000
$ws.Range("L54").Value = 30000
$ws.Range("N54").Value = -31118
$ws.Range("H113").Value = 806.3714
$ws.Range("J113").Value = 821.8889
$ws.Range("L113").Value = 2465.6667
$ws.Range("N113").Value = -6805.6667

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 20035.285
$ws.Range("I102").Value = 22457.834
$ws.Range("K102").Value = 22457.834
$ws.Range("M102").Value = -20835.834
$ws.Range("H113").Value = 6193.5
$ws.Range("I113").Value = 7170.7896
$ws.Range("J113").Value = 2479.8
$ws.Range("K113").Value = 7170.7896
$ws.Range("L113").Value = 2479.8
$ws.Range("M113").Value = -5000.7896
$ws.Range("N113").Value = -6819.8
$ws.Range("H126").Value = 9141.25
$ws.Range("I126").Value = 18783
$ws.Range("J126").Value = 3949.5386
$ws.Range("K126").Value = 56349
$ws.Range("L126").Value = 11848.6158
$ws.Range("M126").Value = -53879
$ws.Range("N126").Value = -16788.6158
$ws.Range("H132").Value = 2133.9312
$ws.Range("I132").Value = 2210.9614
$ws.Range("K132").Value = 6632.8842
$ws.Range("M132").Value = -4102.8842

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1198.8948
$ws.Range("I55").Value = 481.58334
$ws.Range("K55").Value = 481.58334
$ws.Range("M55").Value = -308.58334
$ws.Range("H74").Value = 37199.8
$ws.Range("I74").Value = 28999.75
$ws.Range("K74").Value = 28999.75
$ws.Range("M74").Value = -28001.75
$ws.Range("H77").Value = 37199.8
$ws.Range("I77").Value = 28999.75
$ws.Range("K77").Value = 86999.25
$ws.Range("M77").Value = -82007.25
$ws.Range("H129").Value = 69214
$ws.Range("J129").Value = 69214
$ws.Range("L129").Value = 69214
$ws.Range("N129").Value = -79214
$ws.Range("H132").Value = 393981.47
$ws.Range("I132").Value = 533485.0600000001
$ws.Range("J132").Value = 3371.5
$ws.Range("K132").Value = 1600455.18
$ws.Range("L132").Value = 10114.5
$ws.Range("M132").Value = -1597925.18
$ws.Range("N132").Value = -15174.5
$ws.Range("H136").Value = 4718.7383
$ws.Range("I136").Value = 3261.9473
$ws.Range("K136").Value = 9785.841899999999
$ws.Range("M136").Value = -7235.841899999999

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 45461016
$ws.Range("I2").Value = 45461016
$ws.Range("K2").Value = 45461016
$ws.Range("M2").Value = -45460904
$ws.Range("H4").Value = 860.52
$ws.Range("I4").Value = 175.75
$ws.Range("J4").Value = 3599.6
$ws.Range("K4").Value = 175.75
$ws.Range("L4").Value = 3599.6
$ws.Range("M4").Value = -62.75
$ws.Range("N4").Value = -3825.6
$ws.Range("H113").Value = 1855.6154
$ws.Range("I113").Value = 1011.8095
$ws.Range("K113").Value = 3035.4285
$ws.Range("M113").Value = -865.4285
$ws.Range("H118").Value = 74999.5
$ws.Range("J118").Value = 74999.5
$ws.Range("L118").Value = 74999.5
$ws.Range("N118").Value = -78313.5
$ws.Range("H129").Value = 75000
$ws.Range("J129").Value = 75000
$ws.Range("L129").Value = 75000
$ws.Range("N129").Value = -85000
$ws.Range("H132").Value = 6880.5
$ws.Range("I132").Value = 7747.074
$ws.Range("J132").Value = 4540.75
$ws.Range("K132").Value = 23241.222
$ws.Range("L132").Value = 13622.25
$ws.Range("M132").Value = -20711.222
$ws.Range("N132").Value = -18682.25
$ws.Range("H136").Value = 386684.3
$ws.Range("I136").Value = 441536.22
$ws.Range("K136").Value = 1324608.66
$ws.Range("M136").Value = -1322058.66

